# Apply targeted value corrections to Sheet1 ("TABLO 1")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("B8").Value = 0.1
$ws.Range("G8").Value = 0.2200231481481482

# Row 9
$ws.Range("E9").Value = 0.1
$ws.Range("G9").Value = 0.2200231481481482

# Row 10
$ws.Range("C10").Value = 0.7
$ws.Range("G10").Value = 0.7400077160493828

# Row 12
$ws.Range("E12").Value = 0.3
$ws.Range("G12").Value = 0.4600180041152263
